$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing existing rows 4-5 down to 5-6
$ws.Rows.Item(4).Insert()

# Fill in the new row 4 with the new review data
$ws.Cells.Item(4, 1).Value = 5
$ws.Cells.Item(4, 2).Value = "Agradecido pelo apoio!!!"
$ws.Cells.Item(4, 3).Value = 45895.79028740741
$ws.Cells.Item(4, 4).Value = "ZDc1MjI5MGQtODAzYy00Y2EzLThlYTktY2ZkOGY5ZmJlNDI5OjU3MDE2"

# Match the style (number format) used by the other date cells in column C
$ws.Cells.Item(4, 3).NumberFormat = $ws.Cells.Item(5, 3).NumberFormat
